# PROS-9738 - CCRU - New POS 2019 KPIs
#
# The "FF" sheet's filter/header setup is being corrected: the AutoFilter
# (and its backing _FilterDatabase defined name) should cover only the
# header row A1:AM1 instead of the whole data block A1:AM52, a fresh
# _FilterDatabase_0_0_0_0 bookmark is recorded for the new filter state,
# and the "Activation, Other" tag that had drifted one column to the right
# (into W) is moved back into its proper V column for every affected KPI
# row. The view is also reset back to the top-left of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Defined names: swap which _xlnm._FilterDatabase entry is the
#        hidden "whole table" one vs. the visible "header only" one, then
#        record the new incremental filter-database bookmark. ----------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "FF!_FilterDatabase") {
        if ($n.Visible -eq $false) {
            $hiddenName = $n
        } else {
            $visibleName = $n
        }
    }
}
$hiddenName.RefersTo = "=FF!`$A`$1:`$AM`$1"
$visibleName.RefersTo = "=FF!`$A`$1:`$AM`$52"

$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0", "=FF!`$A`$1:`$AM`$1")

# --- 2. AutoFilter: re-apply over the header row only (A1:AM1) instead of
#        the full A1:AM52 block. ---------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AM1").AutoFilter()

# --- 3. Move the misplaced "Activation, Other" tag from column W back to
#        column V for each affected row. --------------------------------
$rows = @(4,5,6,7,8,9,10,11,12,13,14,15,16,18,19,21,22,24,25,26,28,29,30,31)
foreach ($r in $rows) {
    $wCell = $ws.Cells.Item($r, 23)
    $vCell = $ws.Cells.Item($r, 22)
    $vCell.Value2 = $wCell.Value2
    $wCell.Value2 = ""
}

# --- 4. Reset the view back to the top-left corner and park the active
#        selection on W7. ------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("W7").Select()
